$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = New-Object 'object[,]' 24,8

$data[0,0] = 29.12242425684365
$data[0,1] = 23.02353161545612
$data[0,2] = 5.457562313420514
$data[0,3] = 29.30942974478529
$data[0,4] = 42.83252900737337
$data[0,5] = 2.069037592771776
$data[0,6] = 3.219152413621614
$data[0,7] = 3.413434991520969
$data[1,0] = 27.1374243989585
$data[1,1] = 21.40733941091886
$data[1,2] = 5.337267986424786
$data[1,3] = 27.25656006138376
$data[1,4] = 40.45473511116975
$data[1,5] = 2.078198967281577
$data[1,6] = 2.859703050369023
$data[1,7] = 3.074111969351026
$data[2,0] = 25.87105930567509
$data[2,1] = 20.36552448339359
$data[2,2] = 5.260392023805906
$data[2,3] = 25.9310335883985
$data[2,4] = 38.9386664220199
$data[2,5] = 2.083948925629955
$data[2,6] = 2.633168838164859
$data[2,7] = 2.861811536839314
$data[3,0] = 25.32578769390431
$data[3,1] = 19.93426239784606
$data[3,2] = 5.224182377235195
$data[3,3] = 25.37362898793707
$data[3,4] = 38.28229444321997
$data[3,5] = 2.086343736601319
$data[3,6] = 2.538891211510318
$data[3,7] = 2.774232242323455
$data[4,0] = 25.21926499381419
$data[4,1] = 19.86909335459384
$data[4,2] = 5.213083032900629
$data[4,3] = 25.27925579919151
$data[4,4] = 38.14214131045867
$data[4,5] = 2.086765773327855
$data[4,6] = 2.522620236960007
$data[4,7] = 2.759813657969449
$data[5,0] = 25.82386594742595
$data[5,1] = 20.37891803677698
$data[5,2] = 5.24621841167165
$data[5,3] = 25.92137325416789
$data[5,4] = 38.84853422600711
$data[5,5] = 2.084041586545941
$data[5,6] = 2.630487630487425
$data[5,7] = 2.860669597260886
$data[6,0] = 28.39101441007527
$data[6,1] = 22.49958675012746
$data[6,2] = 5.398999121240402
$data[6,3] = 28.61202614888293
$data[6,4] = 41.92355994330617
$data[6,5] = 2.072250082095517
$data[6,6] = 3.094540212180184
$data[6,7] = 3.296299062931547
$data[7,0] = 33.0895634734422
$data[7,1] = 26.24400022884295
$data[7,2] = 5.698493296056921
$data[7,3] = 33.386799994948
$data[7,4] = 47.63182860177298
$data[7,5] = 2.049914329690132
$data[7,6] = 3.964171270026597
$data[7,7] = 4.12457982399691
$data[8,0] = 36.0434723756846
$data[8,1] = 28.63441750164197
$data[8,2] = 5.830995088903073
$data[8,3] = 35.71059044858454
$data[8,4] = 51.13791178948935
$data[8,5] = 2.034521241778429
$data[8,6] = 4.534438610626464
$data[8,7] = 4.698664340282644
$data[9,0] = 36.02269999258484
$data[9,1] = 28.45708531374927
$data[9,2] = 5.314150617902084
$data[9,3] = 29.16197763703897
$data[9,4] = 49.38345475826942
$data[9,5] = 2.03254108020656
$data[9,6] = 4.85817995592064
$data[9,7] = 4.774015261544365
$data[10,0] = 35.45132614581257
$data[10,1] = 27.80348711356223
$data[10,2] = 4.902634528648329
$data[10,3] = 23.22781200176598
$data[10,4] = 47.2578856750118
$data[10,5] = 2.033596121678061
$data[10,6] = 5.582613299285839
$data[10,7] = 4.732397001373103
$data[11,0] = 34.36056989042019
$data[11,1] = 26.74867764636883
$data[11,2] = 4.535521992404589
$data[11,3] = 17.30423646035947
$data[11,4] = 44.55714983913662
$data[11,5] = 2.037103107177901
$data[11,6] = 6.5172900259235
$data[11,7] = 4.598145302969217
$data[12,0] = 33.34101321820851
$data[12,1] = 25.82294158201565
$data[12,2] = 4.309114174977028
$data[12,3] = 13.23754183903135
$data[12,4] = 42.3561915673402
$data[12,5] = 2.040608279471432
$data[12,6] = 7.264252431868699
$data[12,7] = 4.464898670560597
$data[13,0] = 32.9648862501232
$data[13,1] = 25.51230040168827
$data[13,2] = 4.257088291712681
$data[13,3] = 12.24706028960986
$data[13,4] = 41.66581721299364
$data[13,5] = 2.042053549585688
$data[13,6] = 7.433725654074647
$data[13,7] = 4.411213507189735
$data[14,0] = 31.87608673764142
$data[14,1] = 24.67704261173039
$data[14,2] = 4.270210178006453
$data[14,3] = 11.94105676579287
$data[14,4] = 40.52409631443208
$data[14,5] = 2.047802731150009
$data[14,6] = 7.12150906140563
$data[14,7] = 4.199051954891138
$data[15,0] = 31.59081828550146
$data[15,1] = 24.54023630133161
$data[15,2] = 4.392134655494539
$data[15,3] = 13.92327043899569
$data[15,4] = 40.8602284027189
$data[15,5] = 2.050415076146294
$data[15,6] = 6.426382675262424
$data[15,7] = 4.104955505081555
$data[16,0] = 32.01266221312296
$data[16,1] = 24.99574270966277
$data[16,2] = 4.659682672957626
$data[16,3] = 18.43728141267208
$data[16,4] = 42.57297928984077
$data[16,5] = 2.050369575357616
$data[16,6] = 5.417158518650544
$data[16,7] = 4.10792768057442
$data[17,0] = 32.91428265727141
$data[17,1] = 25.92709831785335
$data[17,2] = 5.046921090650943
$data[17,3] = 24.73815109591906
$data[17,4] = 45.11864271034842
$data[17,5] = 2.04795465780415
$data[17,6] = 4.505350713937855
$data[17,7] = 4.201341048233777
$data[18,0] = 35.18532413088438
$data[18,1] = 28.05782744234055
$data[18,2] = 5.75501748214636
$data[18,3] = 35.0635044301351
$data[18,4] = 50.01271585242392
$data[18,5] = 2.038744722034453
$data[18,6] = 4.378823149754055
$data[18,7] = 4.545497348120654
$data[19,0] = 37.60410681397683
$data[19,1] = 30.04120099125489
$data[19,2] = 5.965210850010251
$data[19,3] = 38.14978486046395
$data[19,4] = 53.24885792036237
$data[19,5] = 2.026048842488905
$data[19,6] = 4.886079521493445
$data[19,7] = 5.01815482679849
$data[20,0] = 39.08938018609194
$data[20,1] = 31.20984196780294
$data[20,2] = 6.081779563242853
$data[20,3] = 39.6612641198891
$data[20,4] = 55.20015203280625
$data[20,5] = 2.01799549018383
$data[20,6] = 5.192950727496313
$data[20,7] = 5.317053975957242
$data[21,0] = 38.3363028547522
$data[21,1] = 30.57376622537857
$data[21,2] = 6.034294527975277
$data[21,3] = 38.86161477420278
$data[21,4] = 54.2361464001029
$data[21,5] = 2.022226197195321
$data[21,6] = 5.030911046414971
$data[21,7] = 5.158508503327936
$data[22,0] = 35.31861004282749
$data[22,1] = 28.11630933623611
$data[22,2] = 5.827425069724348
$data[22,3] = 35.73341028950765
$data[22,4] = 50.37637003397209
$data[22,5] = 2.038380818390848
$data[22,6] = 4.411048834133488
$data[22,7] = 4.555896700817713
$data[23,0] = 31.81816931407905
$data[23,1] = 25.29905695135953
$data[23,2] = 5.596600415070315
$data[23,3] = 32.14356109270349
$data[23,4] = 46.01060943009153
$data[23,5] = 2.055980476337377
$data[23,6] = 3.730515494968214
$data[23,7] = 3.901720770525965

$ws.Range("B2:I25").Value = $data

$pdata = New-Object 'object[,]' 24,1
$pdata[0,0] = 13.37876734696289
$pdata[1,0] = 13.44169929310584
$pdata[2,0] = 13.48005957020183
$pdata[3,0] = 13.49076358101514
$pdata[4,0] = 13.48663518098441
$pdata[5,0] = 13.46411645950251
$pdata[6,0] = 13.3793457196131
$pdata[7,0] = 13.24067711716932
$pdata[8,0] = 13.06500675699586
$pdata[9,0] = 12.41367954742661
$pdata[10,0] = 11.97410097664093
$pdata[11,0] = 11.65383825406395
$pdata[12,0] = 11.49213779797367
$pdata[13,0] = 11.46914446365359
$pdata[14,0] = 11.6109875496389
$pdata[15,0] = 11.7975146447877
$pdata[16,0] = 12.08228754755286
$pdata[17,0] = 12.44996205798263
$pdata[18,0] = 13.05914361158725
$pdata[19,0] = 13.03976532144044
$pdata[20,0] = 13.01478655940469
$pdata[21,0] = 13.04675674622192
$pdata[22,0] = 13.14098249754714
$pdata[23,0] = 13.2493814597989

$ws.Range("P2:P25").Value = $pdata
